# Add a new worksheet "N=200000" at the end of the workbook, matching the
# layout/content used by the other "shell sort details" sheets.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "N=200000"

# Header row
$ws.Cells.Item(1, 1).Value = "Execução"
$ws.Cells.Item(1, 2).Value = "Tempo (ms)"

# Individual run timings
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "575.1212 ms"

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "578.8379 ms"

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "561.8901 ms"

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "589.7360 ms"

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "613.6458 ms"

# Summary rows
$ws.Cells.Item(7, 1).Value = "Média"
$ws.Cells.Item(7, 2).Value = "583.8462 ms"

$ws.Cells.Item(8, 1).Value = "Desvio Padrão"
$ws.Cells.Item(8, 2).Value = "19.4035 ms"

$ws.Range("A1").Select()
